$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold numeric-looking text values
# (e.g. "0.9979", "  +5.45%  ") in the source data. Force text format
# first so Excel does not silently convert them to numbers, then
# restore the default "Normal" style so no stray cell styles remain.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.037.21"
$ws.Range("E2").Value = "  +5.45%  "
$ws.Range("D3").Value = "1.780.26"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "242.94"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "0.9980"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "0.4885"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").Value = "0.2656"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "0.06238"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "1.780.22"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").Value = "16.34"
$ws.Range("E11").Value = "  +4.31%  "
$ws.Range("D12").Value = "0.06995"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "0.6176"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "4.616"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").Value = "79.38"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "27.977.98"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("D17").Value = "0.9988"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "0.9979"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "0.000007207"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").Value = "11.79"
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("D21").Value = "2.015.07"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "4.578"
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("D23").Value = "8.668"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "5.185"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "141.79"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Value = "1.870"
$ws.Range("E27").Value = "  +7.23%  "
$ws.Range("D28").Value = "109.26"
$ws.Range("D29").Value = "1.391"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("D30").Value = "4.091"
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").Value = "0.08270"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("D33").Value = "0.04728"
$ws.Range("E33").Value = "  +4.67%  "
$ws.Range("D34").Value = "1.062"
$ws.Range("E34").Value = "  +6.44%  "
$ws.Range("D35").Value = "2.599"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").Value = "0.6357"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").Value = "0.9410"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "2.594"
$ws.Range("E38").Value = "  +8.08%  "
$ws.Range("D39").Value = "2.060"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "5.845"
$ws.Range("E40").Value = "  +5.91%  "
$ws.Range("D41").Value = "0.01534"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").Value = "0.9987"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "100.05"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "0.3932"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "7.142"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("D46").Value = "0.1195"
$ws.Range("E46").Value = "  +3.36%  "
$ws.Range("D47").Value = "0.05402"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "7.916"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "30.40"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "1.274"
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("D51").Value = "52.41"
$ws.Range("E51").Value = "  +1.88%  "

$priceRange.Style = "Normal"
